$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.508.19'
$ws.Range('E2').Value = '  -2.75%  '
$ws.Range('D3').Value = '1.656.30'
$ws.Range('E3').Value = '  -4.15%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''213.89'
$ws.Range('E5').Value = '  -2.50%  '
$ws.Range('E6').Value = '  -2.51%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''24.02'
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('E9').Value = '  -2.30%  '
$ws.Range('D10').Value = '''0.0617'
$ws.Range('E10').Value = '  -3.16%  '
$ws.Range('D11').Value = '''0.0877'
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').Value = '1.887.87'
$ws.Range('E12').Value = '  -4.29%  '
$ws.Range('D13').Value = '1.655.04'
$ws.Range('E13').Value = '  -4.22%  '
$ws.Range('D14').Value = '''4.13'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').Value = '''0.560'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').Value = '''65.65'
$ws.Range('E16').Value = '  -3.00%  '
$ws.Range('D17').Value = '27.481.19'
$ws.Range('E17').Value = '  -2.76%  '
$ws.Range('D18').Value = '''239.99'
$ws.Range('E18').Value = '  -2.02%  '
$ws.Range('D19').Value = '0.0₃0727'
$ws.Range('E19').Value = '  -3.38%  '
$ws.Range('D20').Value = '''7.53'
$ws.Range('E20').Value = '  -4.98%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').Value = '''4.44'
$ws.Range('E22').Value = '  -4.28%  '
$ws.Range('D23').Value = '''9.29'
$ws.Range('E23').Value = '  -3.95%  '
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('D25').Value = '''145.33'
$ws.Range('E25').Value = '  -2.71%  '
$ws.Range('D26').Value = '''7.18'
$ws.Range('E26').Value = '  -4.20%  '
$ws.Range('D27').Value = '''16.18'
$ws.Range('E27').Value = '  -2.84%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  -2.34%  '
$ws.Range('E30').Value = '  -0.21%  '
$ws.Range('D31').Value = '''0.0498'
$ws.Range('E31').Value = '  -3.61%  '
$ws.Range('D32').Value = '''3.31'
$ws.Range('E32').Value = '  -3.41%  '
$ws.Range('D33').Value = '1.441.15'
$ws.Range('E33').Value = '  -2.78%  '
$ws.Range('D34').Value = '''3.09'
$ws.Range('E34').Value = '  -5.30%  '
$ws.Range('D35').Value = '''1.55'
$ws.Range('E35').Value = '  -5.56%  '
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('D37').Value = '''0.916'
$ws.Range('E37').Value = '  -6.50%  '
$ws.Range('D38').Value = '''0.0170'
$ws.Range('E38').Value = '  -3.02%  '
$ws.Range('D39').Value = '''0.568'
$ws.Range('E39').Value = '  -6.11%  '
$ws.Range('E40').Value = '  -2.93%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '''66.23'
$ws.Range('E42').Value = '  -4.99%  '
$ws.Range('E43').Value = '  -4.01%  '
$ws.Range('D44').Value = '''0.794'
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('D45').Value = '''2.21'
$ws.Range('E45').Value = '  -3.39%  '
$ws.Range('D46').Value = '1.797.55'
$ws.Range('D47').Value = '''1.69'
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('D48').Value = '''88.29'
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('E49').Value = '  -5.99%  '
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('D51').Value = '''7.78'
$ws.Range('E51').Value = '  -4.55%  '
